$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 12: plain value changes (no formulas)
$ws.Range("I12").Value = 4188377156.3100057
$ws.Range("J12").Value = 3588029419

# Row 13: plain value changes (no formulas)
$ws.Range("I13").Value = 1012006300.0300001
$ws.Range("J13").Value = 9569343401.6000004

# Row 14: plain value changes (no formulas)
$ws.Range("I14").Value = -44319159.290000051
$ws.Range("J14").Value = 146268235.09999999

# Row 16: plain value changes (no formulas)
$ws.Range("I16").Value = -162861893.56999999
$ws.Range("J16").Value = -193292161.30000001

# Row 18: I18 becomes a formula SUM(I12:I17); J18 keeps its formula, value recalculates
$ws.Range("I18").Formula = "=SUM(I12:I17)"

# Row 19: plain value change (no formula)
$ws.Range("I19").Value = -1160500000.0000002

# Row 21: I21 becomes a formula SUM(I18:I20); J21 keeps its formula, value recalculates
$ws.Range("I21").Formula = "=SUM(I18:I20)"

# Row 26: plain value changes (no formulas)
$ws.Range("I26").Value = 1010658958.9880759
$ws.Range("J26").Value = 1031977291

# Force recalculation of all formulas so dependent cells (I23, J23, I25, J25, J28, I29) update
$excel.CalculateFullRebuild()
